$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# --- Sheet LP1912 ---
$ws1.Cells.Item(2,1).Value = "Última actualización: 11:13:15"
$ws1.Cells.Item(3,1).Value = "Total filas: 129"
$ws1.Cells.Item(20,1).Value = "05:57:04"
$ws1.Cells.Item(20,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(20,4).Value = 84
$ws1.Cells.Item(21,1).Value = "06:16:41"
$ws1.Cells.Item(21,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(21,4).Value = 65
$ws1.Cells.Item(86,1).Value = "10:04:30"
$ws1.Cells.Item(86,3).Value = "14_ABASTO"
$ws1.Cells.Item(86,4).Value = 25
$ws1.Cells.Item(87,1).Value = "08:38:24"
$ws1.Cells.Item(87,3).Value = "15_ABASTO"
$ws1.Cells.Item(87,4).Value = 111
$ws1.Cells.Item(102,1).Value = "10:56:15"
$ws1.Cells.Item(102,3).Value = "14_ABASTO"
$ws1.Cells.Item(102,4).Value = 19
$ws1.Cells.Item(103,1).Value = "09:22:34"
$ws1.Cells.Item(103,3).Value = "15X38_ABASTO"
$ws1.Cells.Item(103,4).Value = 113
$ws1.Cells.Item(104,1).Value = "11:13:15"
$ws1.Cells.Item(104,2).Value = "11:17"
$ws1.Cells.Item(104,3).Value = "14_ABASTO"
$ws1.Cells.Item(104,4).Value = 4
$ws1.Cells.Item(105,1).Value = "10:36:50"
$ws1.Cells.Item(105,2).Value = "11:25"
$ws1.Cells.Item(105,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(105,4).Value = 49
$ws1.Cells.Item(107,1).Value = "10:04:30"
$ws1.Cells.Item(107,2).Value = "11:29"
$ws1.Cells.Item(107,3).Value = "10_OLMOS"
$ws1.Cells.Item(107,4).Value = 85
$ws1.Cells.Item(108,1).Value = "10:36:50"
$ws1.Cells.Item(108,2).Value = "11:30"
$ws1.Cells.Item(108,4).Value = 54
$ws1.Cells.Item(109,2).Value = "11:31"
$ws1.Cells.Item(109,3).Value = "215C_EL PATO"
$ws1.Cells.Item(109,4).Value = 87
$ws1.Cells.Item(110,1).Value = "10:04:30"
$ws1.Cells.Item(110,2).Value = "11:41"
$ws1.Cells.Item(110,4).Value = 97
$ws1.Cells.Item(111,1).Value = "10:56:15"
$ws1.Cells.Item(111,2).Value = "11:42"
$ws1.Cells.Item(111,3).Value = "215B_EL PATO"
$ws1.Cells.Item(111,4).Value = 46
$ws1.Cells.Item(112,1).Value = "10:04:30"
$ws1.Cells.Item(112,2).Value = "11:45"
$ws1.Cells.Item(112,3).Value = "15X38_ABASTO"
$ws1.Cells.Item(112,4).Value = 101
$ws1.Cells.Item(113,1).Value = "10:56:15"
$ws1.Cells.Item(113,2).Value = "11:46"
$ws1.Cells.Item(113,4).Value = 50
$ws1.Cells.Item(114,1).Value = "10:49:38"
$ws1.Cells.Item(114,2).Value = "11:47"
$ws1.Cells.Item(114,4).Value = 58
$ws1.Cells.Item(115,2).Value = "11:48"
$ws1.Cells.Item(115,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(115,4).Value = 72
$ws1.Cells.Item(116,1).Value = "11:13:15"
$ws1.Cells.Item(116,2).Value = "11:51"
$ws1.Cells.Item(116,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(116,4).Value = 38
$ws1.Cells.Item(117,1).Value = "10:36:50"
$ws1.Cells.Item(117,2).Value = "11:52"
$ws1.Cells.Item(117,3).Value = "225_GOMEZ"
$ws1.Cells.Item(117,4).Value = 76
$ws1.Cells.Item(118,1).Value = "10:04:30"
$ws1.Cells.Item(118,2).Value = "11:53"
$ws1.Cells.Item(118,3).Value = "225_GOMEZ"
$ws1.Cells.Item(118,4).Value = 109
$ws1.Cells.Item(119,1).Value = "10:04:30"
$ws1.Cells.Item(119,2).Value = "11:58"
$ws1.Cells.Item(119,3).Value = "17_ROMERO"
$ws1.Cells.Item(119,4).Value = 114
$ws1.Cells.Item(120,2).Value = "12:05"
$ws1.Cells.Item(120,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(120,4).Value = 89
$ws1.Cells.Item(121,1).Value = "10:56:15"
$ws1.Cells.Item(121,2).Value = "12:06"
$ws1.Cells.Item(121,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(121,4).Value = 70
$ws1.Cells.Item(122,2).Value = "12:10"
$ws1.Cells.Item(122,3).Value = "15_ABASTO"
$ws1.Cells.Item(122,4).Value = 94
$ws1.Cells.Item(123,1).Value = "10:36:50"
$ws1.Cells.Item(123,2).Value = "12:10"
$ws1.Cells.Item(123,3).Value = "16_P MOR-SANTA ANA"
$ws1.Cells.Item(123,4).Value = 94
$ws1.Cells.Item(124,1).Value = "11:13:15"
$ws1.Cells.Item(124,2).Value = "12:17"
$ws1.Cells.Item(124,3).Value = "10_OLMOS"
$ws1.Cells.Item(124,4).Value = 64
$ws1.Cells.Item(125,1).Value = "10:36:50"
$ws1.Cells.Item(125,2).Value = "12:21"
$ws1.Cells.Item(125,3).Value = "215C_EL PATO"
$ws1.Cells.Item(125,4).Value = 105
$ws1.Cells.Item(126,2).Value = "12:22"
$ws1.Cells.Item(126,3).Value = "215C_EL PATO"
$ws1.Cells.Item(126,4).Value = 86
$ws1.Cells.Item(127,1).Value = "11:13:15"
$ws1.Cells.Item(127,2).Value = "12:29"
$ws1.Cells.Item(127,3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(127,4).Value = 76
$ws1.Cells.Item(128,1).Value = "10:36:50"
$ws1.Cells.Item(128,2).Value = "12:32"
$ws1.Cells.Item(128,3).Value = "14_ABASTO"
$ws1.Cells.Item(128,4).Value = 116
$ws1.Cells.Item(129,1).Value = "10:56:15"
$ws1.Cells.Item(129,2).Value = "12:33"
$ws1.Cells.Item(129,3).Value = "14_ABASTO"
$ws1.Cells.Item(129,4).Value = 97
$ws1.Cells.Item(130,1).Value = "10:56:15"
$ws1.Cells.Item(130,2).Value = "12:33"
$ws1.Cells.Item(130,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(130,4).Value = 97
$ws1.Cells.Item(130,5).Value = "LP1912"
$ws1.Cells.Item(131,1).Value = "10:36:50"
$ws1.Cells.Item(131,2).Value = "12:34"
$ws1.Cells.Item(131,3).Value = "15_ABASTO"
$ws1.Cells.Item(131,4).Value = 118
$ws1.Cells.Item(131,5).Value = "LP1912"
$ws1.Cells.Item(132,1).Value = "10:49:38"
$ws1.Cells.Item(132,2).Value = "12:36"
$ws1.Cells.Item(132,3).Value = "27_EL RETIRO"
$ws1.Cells.Item(132,4).Value = 107
$ws1.Cells.Item(132,5).Value = "LP1912"
$ws1.Cells.Item(133,1).Value = "10:49:38"
$ws1.Cells.Item(133,2).Value = "12:48"
$ws1.Cells.Item(133,3).Value = "16_SANTA ANA"
$ws1.Cells.Item(133,4).Value = 119
$ws1.Cells.Item(133,5).Value = "LP1912"
$ws1.Cells.Item(134,1).Value = "11:13:15"
$ws1.Cells.Item(134,2).Value = "13:03"
$ws1.Cells.Item(134,3).Value = "11_ETCHEVERRY"
$ws1.Cells.Item(134,4).Value = 110
$ws1.Cells.Item(134,5).Value = "LP1912"

# --- Sheet LP1912-215 ---
$ws2.Cells.Item(2,1).Value = "Última actualización: 11:13:15"

# --- Sheet 6203-6173 ---
$ws3.Cells.Item(2,1).Value = "Última actualización: 11:13:15"
$ws3.Cells.Item(3,1).Value = "Total filas: 20"
$ws3.Cells.Item(25,1).Value = "11:13:15"
$ws3.Cells.Item(25,2).Value = "13:12"
$ws3.Cells.Item(25,3).Value = "215C_LA PLATA"
$ws3.Cells.Item(25,4).Value = 119
$ws3.Cells.Item(25,5).Value = "L6203"
